$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.004.95'
$ws.Range('E2').Value = '  -0.03%  '
$ws.Range('D3').Value = '2.299.16'
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '300.52'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.10%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '99.56'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.37%  '
$ws.Range('E7').Value = '  -0.61%  '
$ws.Range('E9').Value = '  +1.78%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.21'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +7.75%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0790'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.46%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.116'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.84%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '18.03'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +5.57%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.91'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.99%  '
$ws.Range('D15').Value = '2.660.52'
$ws.Range('E15').Value = '  +0.03%  '
$ws.Range('D16').Value = '2.374.88'
$ws.Range('E16').Value = '  +3.77%  '
$ws.Range('E17').Value = '  -1.25%  '
$ws.Range('D18').Value = '42.913.40'
$ws.Range('E18').Value = '  -0.18%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.69'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +9.53%  '
$ws.Range('D20').Value = '0.0₃0904'
$ws.Range('E20').Value = '  +0.32%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.12'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.97%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '67.86'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.08%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '235.61'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.44%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.18'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +7.37%  '
$ws.Range('E26').Value = '  -0.37%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '24.90'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.74%  '
$ws.Range('E28').Value = '  +14.60%  '
$ws.Range('B29').Value = 'Monero'
$ws.Range('C29').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '167.88'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.96%  '
$ws.Range('B30').Value = 'InjectiveProtocol'
$ws.Range('C30').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '34.51'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.71%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '9.11'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.19%  '
$ws.Range('E32').Value = '  -0.13%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.02'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.49%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '17.57'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +4.47%  '
$ws.Range('E35').Value = '  -1.44%  '
$ws.Range('E36').Value = '  +0.96%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0688'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.99%  '
$ws.Range('B38').Value = 'LidoDAOToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.81'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.64%  '
$ws.Range('B39').Value = 'ARBITRUM'
$ws.Range('C39').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.78'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.93%  '
$ws.Range('E40').Value = '  -0.46%  '
$ws.Range('E41').Value = '  +0.18%  '
$ws.Range('E42').Value = '  -4.91%  '
$ws.Range('E43').Value = '  +3.31%  '
$ws.Range('D44').Value = '1.975.30'
$ws.Range('E44').Value = '  -0.34%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '10.14'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.68%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.90'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.82%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '17.38'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.79%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '55.21'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +3.50%  '
$ws.Range('E49').Value = '  +4.01%  '
$ws.Range('D50').Value = '2.523.04'
$ws.Range('E50').Value = '  -0.11%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '70.59'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.83%  '
